$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell C5 holding a new shared string; this extends the used range
# to A1:C5 and becomes the active selection (matches the target diff).
$ws.Range("C5").Value = "fkeop34opf4f"
$ws.Range("C5").Select()
